$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Khan"
$ws.Range("B3").Value = "yusuf"
$ws.Range("C3").Value = "res"

$ws.Range("B3").Select()
